function Set-TextValue($range, $value) {
    $escaped = $value -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.219.90'
$ws.Range("E2").Value = '  -1.45%  '

$ws.Range("D3").Value = '2.469.01'
$ws.Range("E3").Value = '  -1.91%  '

$ws.Range("E4").Value = '  +0.04%  '

Set-TextValue $ws.Range("D5") '519.13'
$ws.Range("E5").Value = '  -2.83%  '

$ws.Range("E6").Value = '  -1.87%  '

$ws.Range("E7").Value = '  +0.07%  '

Set-TextValue $ws.Range("D8") '0.558'
$ws.Range("E8").Value = '  -1.68%  '

$ws.Range("D9").Value = '2.479.28'
$ws.Range("E9").Value = '  -1.82%  '

$ws.Range("E10").Value = '  -3.43%  '

$ws.Range("E11").Value = '  -0.93%  '

Set-TextValue $ws.Range("D12") '5.30'
$ws.Range("E12").Value = '  -2.15%  '

$ws.Range("E13").Value = '  -2.74%  '

$ws.Range("D14").Value = '2.920.83'
$ws.Range("E14").Value = '  -1.29%  '

$ws.Range("D15").Value = '58.137.38'
$ws.Range("E15").Value = '  -1.50%  '

Set-TextValue $ws.Range("D16") '21.94'
$ws.Range("E16").Value = '  -3.88%  '

$ws.Range("E17").Value = '  -2.38%  '

$ws.Range("D18").Value = '2.475.38'
$ws.Range("E18").Value = '  -1.11%  '

Set-TextValue $ws.Range("D19") '10.60'
$ws.Range("E19").Value = '  -4.38%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D20") '4.17'
$ws.Range("E20").Value = '  -2.21%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D21") '319.27'
$ws.Range("E21").Value = '  -1.64%  '

Set-TextValue $ws.Range("D22") '1.00'
$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("E23").Value = '  -3.79%  '

Set-TextValue $ws.Range("D24") '64.62'
$ws.Range("E24").Value = '  -0.75%  '

Set-TextValue $ws.Range("D25") '0.410'
$ws.Range("E25").Value = '  -2.72%  '

$ws.Range("E26").Value = '  -0.28%  '

$ws.Range("E27").Value = '  -1.62%  '

Set-TextValue $ws.Range("D28") '7.34'
$ws.Range("E28").Value = '  -2.85%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D29") '170.18'
$ws.Range("E29").Value = '  +0.95%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0746'
$ws.Range("E30").Value = '  -2.74%  '

Set-TextValue $ws.Range("D31") '6.36'
$ws.Range("E31").Value = '  -1.90%  '

$ws.Range("E32").Value = '  -3.41%  '

$ws.Range("E33").Value = '  +4.33%  '

Set-TextValue $ws.Range("D34") '0.999'
$ws.Range("E34").Value = '  -0.01%  '

Set-TextValue $ws.Range("D35") '0.998'
$ws.Range("E35").Value = '  +0.05%  '

Set-TextValue $ws.Range("D36") '18.03'
$ws.Range("E36").Value = '  -2.10%  '

$ws.Range("E37").Value = '  -4.36%  '

$ws.Range("E38").Value = '  -2.12%  '

Set-TextValue $ws.Range("D39") '36.61'
$ws.Range("E39").Value = '  -0.84%  '

$ws.Range("E40").Value = '  -4.00%  '

Set-TextValue $ws.Range("D41") '0.795'
$ws.Range("E41").Value = '  -1.17%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D42") '5.10'
$ws.Range("E42").Value = '  +0.68%  '

$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D43") '273.71'
$ws.Range("E43").Value = '  -3.07%  '

$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D44") '3.43'
$ws.Range("E44").Value = '  -4.33%  '

$ws.Range("E45").Value = '  -2.05%  '

Set-TextValue $ws.Range("D46") '123.63'
$ws.Range("E46").Value = '  -5.11%  '

Set-TextValue $ws.Range("D47") '0.0907'
$ws.Range("E47").Value = '  -1.77%  '

$ws.Range("E48").Value = '  -2.62%  '

$ws.Range("E49").Value = '  -2.69%  '

Set-TextValue $ws.Range("D50") '16.97'
$ws.Range("E50").Value = '  -2.08%  '

$ws.Range("D51").Value = '1.733.34'
$ws.Range("E51").Value = '  -1.34%  '
